$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8, mirroring the columns of rows 4-7:
# A = FW VER | B = SPIFFS VER | C = Commit | D = Label on gateway | E = Data | F = Note | G = USR
# A8/B8 are entered with a leading apostrophe so Excel stores them as
# text-with-quote-prefix (same as the existing "100"/"101" cells above).
$ws.Range("A8").Value = "'100"
$ws.Range("B8").Value = "'101"
$ws.Range("C8").Value = 967

# Note text goes in before the short "Label on gateway" value so the two
# new shared-string entries land in the same order as the target file.
$ws.Range("F8").Value = "change the send values routine to resolve 1 second spurious alarms (send before check Alarm, High and low variable)/ `nadded OTA GME HTTPS RANGE routine/ added upload log response inside header"
$ws.Range("D8").Value = "V.1.030"

$ws.Range("F8").WrapText = $true

# Date (06/18/2021), formatted the same way as the other Data cells
$ws.Range("E8").Value = 44365
$ws.Range("E8").NumberFormat = "m/d/yy"

# Wrapped note text needs 3 display lines at the default row height
$ws.Rows.Item(8).RowHeight = 43.2

# Move selection like the author's last recorded action
$ws.Range("B9").Select()
